# Apply latest cryptos snapshot (price + 1h volume change, plus
# ranking reshuffle) onto Sheet1, matching the upstream data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '57.990.96'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +0.86%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.125.26'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +0.02%  '
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '530.48'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.81%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '138.54'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.49%  '
$ws.Range('E7').Value = '  -0.10%  '
$ws.Range('B8').Value = 'XRP'
$ws.Range('C8').Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.461'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +4.48%  '
$ws.Range('B9').Value = 'Toncoin'
$ws.Range('C9').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '7.29'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +1.76%  '
$ws.Range('B10').Value = 'Dogecoin'
$ws.Range('C10').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.107'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.67%  '
$ws.Range('B11').Value = 'Cardano'
$ws.Range('C11').Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.408'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +3.38%  '
$ws.Range('B12').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C12').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '3.657.75'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.30%  '
$ws.Range('B13').Value = 'TRON'
$ws.Range('C13').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.136'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +1.38%  '
$ws.Range('B14').Value = 'Avalanche'
$ws.Range('C14').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '25.47'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.81%  '
$ws.Range('B15').Value = 'ShibaInu'
$ws.Range('C15').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.0000163'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.68%  '
$ws.Range('B16').Value = 'WrappedBTC'
$ws.Range('C16').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '58.032.07'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.63%  '
$ws.Range('B17').Value = 'WrappedEther'
$ws.Range('C17').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '3.110.16'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.73%  '
$ws.Range('B18').Value = 'Polkadot'
$ws.Range('C18').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '5.98'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.96%  '
$ws.Range('B19').Value = 'Chainlink'
$ws.Range('C19').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.65'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.80%  '
$ws.Range('B20').Value = 'Uniswap'
$ws.Range('C20').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '8.07'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +2.48%  '
$ws.Range('B21').Value = 'BitcoinCash'
$ws.Range('C21').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '353.26'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.02%  '
$ws.Range('B22').Value = 'Dai'
$ws.Range('C22').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.998'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.17%  '
$ws.Range('B23').Value = 'Litecoin'
$ws.Range('C23').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '68.72'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.40%  '
$ws.Range('B24').Value = 'Polygon'
$ws.Range('C24').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.503'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.58%  '
$ws.Range('B25').Value = 'Kaspa'
$ws.Range('C25').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.167'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -1.56%  '
$ws.Range('B26').Value = 'Binance-PegBSC-USD'
$ws.Range('C26').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.00'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.26%  '
$ws.Range('B27').Value = 'PEPE'
$ws.Range('C27').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.0₃0882'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -4.21%  '
$ws.Range('B28').Value = 'InternetComputer(DFINITY)'
$ws.Range('C28').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.27'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -1.84%  '
$ws.Range('B29').Value = 'RenderToken'
$ws.Range('C29').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '6.12'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -2.67%  '
$ws.Range('B30').Value = 'PancakeSwap'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.86'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -1.50%  '
$ws.Range('B31').Value = 'EthereumClassic'
$ws.Range('C31').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '21.28'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.90%  '
$ws.Range('B32').Value = 'NEARProtocol'
$ws.Range('C32').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '5.00'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +2.54%  '
$ws.Range('B33').Value = 'Fetch.AI'
$ws.Range('C33').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.15'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -2.64%  '
$ws.Range('B34').Value = 'Monero'
$ws.Range('C34').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '158.77'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.96%  '
$ws.Range('B35').Value = 'Aptos'
$ws.Range('C35').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '6.07'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -1.40%  '
$ws.Range('B36').Value = 'EnergySwap'
$ws.Range('C36').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '26.02'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.69%  '
$ws.Range('B37').Value = 'ImmutableX'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.26'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.57%  '
$ws.Range('B38').Value = 'Stacks'
$ws.Range('C38').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.68'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +4.03%  '
$ws.Range('B39').Value = 'Hedera'
$ws.Range('C39').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0670'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.30%  '
$ws.Range('B40').Value = 'Mantle'
$ws.Range('C40').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.699'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.21%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.99'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -2.79%  '
$ws.Range('B42').Value = 'OKB'
$ws.Range('C42').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '37.57'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +2.60%  '
$ws.Range('B43').Value = 'Maker'
$ws.Range('C43').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.391.12'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +2.85%  '
$ws.Range('B44').Value = 'RenzoRestakedETH'
$ws.Range('C44').Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '3.164.34'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.10%  '
$ws.Range('B45').Value = 'FirstDigitalUSD'
$ws.Range('C45').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.00'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.05%  '
$ws.Range('B46').Value = 'VeChain'
$ws.Range('C46').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0267'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -1.30%  '
$ws.Range('B47').Value = 'ONDO'
$ws.Range('C47').Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.980'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.20%  '
$ws.Range('B48').Value = 'Cosmos'
$ws.Range('C48').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '6.03'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.29%  '
$ws.Range('B49').Value = 'InjectiveProtocol'
$ws.Range('C49').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '19.83'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -2.15%  '
$ws.Range('B50').Value = 'SuiNetwork'
$ws.Range('C50').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.737'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -1.94%  '
$ws.Range('B51').Value = 'Stellar'
$ws.Range('C51').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0908'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +2.36%  '
